# Lots of refactoring, removed JSON support.
#
# Fills in the student identification fields and the "Remarks" column (M)
# for several criteria rows on the "Evaluation form" sheet, replaces the
# wording of the walking-mode remark, and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Evaluation form")

# --- Student identification (top of the form) ---
$ws.Range("E4").Value = "s1129160"
$ws.Range("E5").Value = "Lucas Ouwens"

# --- Techniques section: per-criterion remarks (column M) ---
$ws.Range("M19").Value = "Skybox, floor"

$ws.Range("M20").Value = "tile mesh / floor model, skybox consists of multiple primitive meshes (6 of them)"
$ws.Rows.Item(20).RowHeight = 45

$ws.Range("M22").Value = "Skybox texture, floor grass texture, floor dirt texture"

$ws.Range("M24").Value = "three types of shaders (includes the fragment shading aspect): Basic, Lambert and Phong."

$ws.Range("M25").Value = "Implemented walking mode, uses mouse movement for the camera (so no I/J/K/L) fallthrough is impossible."

$ws.Range("M29").Value = "Application of strategy pattern for movement controllers, code is documented."

# --- Restore the active selection/view ---
$ws.Range("K30").Select()
